# Add BPSK I/Q display: update loop filter cutoff and gain settings.
$wb = $excel.ActiveWorkbook

$wsLoop = $wb.Worksheets.Item("LoopFilter LPF")
$wsBranch = $wb.Worksheets.Item("Branch LPF")
$wsGains = $wb.Worksheets.Item("Gains")

# LoopFilter LPF: cutoff freq 100 -> 50, Gain 8 -> 4
$wsLoop.Range("B2").Value = 50
$wsLoop.Range("B16").Value = 4

# Update selections to match author's final cursor position
$wsLoop.Activate() | Out-Null
$wsLoop.Range("D18").Select() | Out-Null

$wsBranch.Activate() | Out-Null
$wsBranch.Range("C20").Select() | Out-Null

$wsLoop.Activate() | Out-Null
